$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.529.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.791.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.788.62'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.488'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.418.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.789.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.620.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '511.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000138'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("E28").Value = '  -3.98%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.98%  '
$ws.Range("E32").Value = '  +4.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("E38").Value = '  +7.71%  '
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '459.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.37%  '
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("E43").Value = '  +7.88%  '
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.964.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("E51").Value = '  +1.05%  '
